$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "71.813.06"
$ws.Range("E2").Value = "  +3.65%  "
$ws.Range("D3").Value = "3.684.45"
$ws.Range("E3").Value = "  +8.72%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.06"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.09"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "3.677.10"
$ws.Range("E7").Value = "  +8.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.624"
$ws.Range("E8").Value = "  +5.41%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.201"
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.614"
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.92"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "4.281.17"
$ws.Range("E14").Value = "  +8.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "682.31"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.98"
$ws.Range("E16").Value = "  +4.66%  "
$ws.Range("D17").Value = "71.930.72"
$ws.Range("E17").Value = "  +3.65%  "
$ws.Range("D18").Value = "3.675.88"
$ws.Range("E18").Value = "  +8.50%  "
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.03"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.939"
$ws.Range("E22").Value = "  +3.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.18"
$ws.Range("E23").Value = "  +15.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.80"
$ws.Range("E24").Value = "  +3.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.49"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.03"
$ws.Range("E26").Value = "  +3.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.84"
$ws.Range("E27").Value = "  +5.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.44"
$ws.Range("E29").Value = "  +5.98%  "
$ws.Range("E30").Value = "  +5.31%  "
$ws.Range("E31").Value = "  +6.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.19"
$ws.Range("E32").Value = "  +9.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "574.87"
$ws.Range("E33").Value = "  +4.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.30"
$ws.Range("E34").Value = "  +2.42%  "
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.66"
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("D37").Value = "3.755.02"
$ws.Range("E37").Value = "  +4.22%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +3.47%  "
$ws.Range("D40").Value = "0.0₃0774"
$ws.Range("E40").Value = "  +4.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.40"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("E43").Value = "  +8.90%  "
$ws.Range("E44").Value = "  +2.86%  "
$ws.Range("E45").Value = "  +4.11%  "
$ws.Range("E46").Value = "  +8.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.37"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("E48").Value = "  +4.07%  "
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.10"
$ws.Range("E51").Value = "  +3.13%  "
